$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire data row that contains "leniyadoniv@gmail.com" /
# "lets rewrite the rules of bitcoin" (original row 20). This shifts all
# subsequent rows up by one.
$ws.Rows.Item(20).Delete()

# After the deletion, the review marked "yes" for the review on what is now
# row 19 (echale484@gmail.com / never stop learning. Good guide) and the
# review on what is now row 20 (rotemzinger3@gmail.com / show me the money
# please, formerly row 21) both need to be updated to "no".
$ws.Range("G19").Value = "no"
$ws.Range("G20").Value = "no"

# Update the active selection to match the edited workbook.
$ws.Range("G21").Select()

# Deleting the row does not re-anchor the worksheet's mailto: hyperlinks, so
# rebuild the hyperlink collection from scratch to match the new layout.
$ws.Hyperlinks.Delete()

$links = @(
    @{Cell="C7";  Addr="mailto:jorjkluni03@gmail.com";      Disp="jorjkluni03@gmail.com"},
    @{Cell="D7";  Addr="mailto:vikicrestina@gmail.com";     Disp="vikicrestina@gmail.com"},
    @{Cell="D8";  Addr="mailto:jorjkluni03@gmail.com";      Disp="jorjkluni03@gmail.com"},
    @{Cell="C9";  Addr="mailto:nachumshayil@gmail.com";     Disp="nachumshayil@gmail.com"},
    @{Cell="D9";  Addr="mailto:nachushay@gmail.com";        Disp="nachushay@gmail.com"},
    @{Cell="C10"; Addr="mailto:nevilgreen12@gmail.com";     Disp="nevilgreen12@gmail.com"},
    @{Cell="D10"; Addr="mailto:vikicrestina@gmail.com";     Disp="vikicrestina@gmail.com"},
    @{Cell="C11"; Addr="mailto:snizzvered@gmail.com";       Disp="snizzvered@gmail.com"},
    @{Cell="D11"; Addr="mailto:krigelron@gmail.com";        Disp="krigelron@gmail.com"},
    @{Cell="C12"; Addr="mailto:redvelvetmichael@gmail.com"; Disp="redvelvetmichael@gmail.com"},
    @{Cell="D12"; Addr="mailto:veredsnir12@gmail.com";      Disp="veredsnir12@gmail.com"},
    @{Cell="C13"; Addr="mailto:veredsnir12@gmail.com";      Disp="veredsnir12@gmail.com"},
    @{Cell="D13"; Addr="mailto:kevinkors122@gmail.com";     Disp="kevinkors122@gmail.com"},
    @{Cell="C14"; Addr="mailto:freelancernachus@gmail.com"; Disp="freelancernachus@gmail.com"},
    @{Cell="C15"; Addr="mailto:sm6502345@gmail.com";        Disp="sm6502345@gmail.com"},
    @{Cell="D15"; Addr="mailto:cybworking@gmail.com";       Disp="cybworking@gmail.com"},
    @{Cell="C16"; Addr="mailto:rontiddler560@gmail.com";    Disp="rontiddler560@gmail.com"},
    @{Cell="D16"; Addr="mailto:halachme@gmail.com";         Disp="halachme@gmail.com"},
    @{Cell="D18"; Addr="mailto:itaisenior@gmail.com";       Disp="itaisenior@gmail.com"},
    @{Cell="C20"; Addr="mailto:rotemzinger3@gmail.com";     Disp="rotemzinger3@gmail.com"},
    @{Cell="C21"; Addr="mailto:sinuspai@gmail.com";         Disp="sinuspai@gmail.com"},
    @{Cell="D22"; Addr="mailto:rotemzinger3@gmail.com";     Disp="rotemzinger3@gmail.com"}
)

foreach ($link in $links) {
    $ws.Hyperlinks.Add($ws.Range($link.Cell), $link.Addr, [System.Type]::Missing, [System.Type]::Missing, $link.Disp) | Out-Null
}
